# ImportUsersSampleFile.xlsx — "complete update of code base"
#
# The header row used to mark required columns with a trailing "*"
# (UserName*, Name*, Surname*, EmailAddress*, Password*). That suffix is
# dropped here. The sheet's saved cursor/selection also moved from the
# last populated row (A4:XFD4) to cell F5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Surname"
$ws.Range("D1").Value = "EmailAddress"
$ws.Range("F1").Value = "Password"

# Match the saved view's active cell/selection.
[void]$ws.Range("F5").Select()
